$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for student 190540 ("ليان بنت خالد بن سعد المقذلى"), currently the
# first student of group B2A (row 6), moves down so it becomes the last row
# of what is now the B2D block, with its Group re-labelled to B2E (the block
# immediately below it in the new layout).
#
# Implemented as a pure value shift (no structural row insert/delete) so the
# existing per-row banding styles (which alternate strictly by absolute row
# number, independent of content) are left completely untouched: rows
# 7..223 slide up into 6..222, and the vacated row 222 is filled back in
# with the moved record's data.

$movedId = $ws.Range("A6").Value()
$movedName = $ws.Range("B6").Value()
$movedYear = $ws.Range("C6").Value()
$movedSource = $ws.Range("E6").Value()

# Shift the Name/Year/Group/Source File columns in bulk.
$ws.Range("B6:E222").Value = $ws.Range("B7:E223").Value()

# Shift the Student ID column one row at a time, forcing each value back to
# text (leading apostrophe) so numeric-looking IDs don't get silently
# reinterpreted as numbers by the COM value setter.
for ($r = 6; $r -le 221; $r++) {
    $ws.Cells.Item($r, 1).Value = "'" + $ws.Cells.Item($r + 1, 1).Value()
}

$ws.Range("A222").Value = "'" + $movedId
$ws.Range("B222").Value = $movedName
$ws.Range("C222").Value = $movedYear
$ws.Range("D222").Value = "B2E"
$ws.Range("E222").Value = $movedSource
